$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 41 (pushes existing rows 41-48 down to 42-49)
$ws.Rows.Item(41).Insert()

# Fill the new row 41 with the latest weekly record
$ws.Cells.Item(41, 1).Value = 9
$ws.Cells.Item(41, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(41, 3).Value = "Metropolitana"
$ws.Cells.Item(41, 4).Value = 45015
$ws.Cells.Item(41, 5).Value = 13
$ws.Cells.Item(41, 6).Value = 100112010
$ws.Cells.Item(41, 7).Value = "Achicoria"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 90
$ws.Cells.Item(41, 11).Value = 7000
$ws.Cells.Item(41, 12).Value = 7000
$ws.Cells.Item(41, 13).Value = 7000
$ws.Cells.Item(41, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(41, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(41, 16).Value = 438
$ws.Cells.Item(41, 17).Value = 16
$ws.Cells.Item(41, 18).Value = "Hortaliza"
